# Refresh cryptocurrency price/volume snapshot (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a "Price" (column D) cell while forcing it to be
# stored as literal text. Many of these prices look numeric to Excel's
# input-parser (e.g. "311.61", "1.000") and would otherwise be silently
# reinterpreted as a number (losing trailing zeros / exact formatting), so we
# briefly switch the cell to the Text number format, assign the string, then
# restore the cell's style so no formatting change is left behind.
function Set-PriceText($rangeRef, $text) {
    $cell = $ws.Range($rangeRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "24.889.43"
$ws.Range("E2").Value = "  +1.71%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.709.30"
$ws.Range("E3").Value = "  +1.59%  "

# Row 4 - TetherUSD (price unchanged, volume refreshed)
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
Set-PriceText "D5" "311.61"
$ws.Range("E5").Value = "  +1.27%  "

# Row 6 - USDC
Set-PriceText "D6" "0.9981"
$ws.Range("E6").Value = "  -0.14%  "

# Row 7 - XRP
Set-PriceText "D7" "0.3745"
$ws.Range("E7").Value = "  +0.85%  "

# Row 8 - OKB
Set-PriceText "D8" "49.42"
$ws.Range("E8").Value = "  +3.63%  "

# Row 9 - Cardano
Set-PriceText "D9" "0.3444"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10 - Polygon
Set-PriceText "D10" "1.217"
$ws.Range("E10").Value = "  +3.30%  "

# Row 11 - Dogecoin
Set-PriceText "D11" "0.07537"
$ws.Range("E11").Value = "  +3.59%  "

# Row 12 - BinanceUSD
Set-PriceText "D12" "1.000"
$ws.Range("E12").Value = "  -0.05%  "

# Row 13 - Solana
Set-PriceText "D13" "21.26"
$ws.Range("E13").Value = "  +4.24%  "

# Row 14 - Polkadot
Set-PriceText "D14" "6.310"
$ws.Range("E14").Value = "  +2.43%  "

# Row 15 - Chainlink
Set-PriceText "D15" "7.095"
$ws.Range("E15").Value = "  +4.92%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.704.91"
$ws.Range("E16").Value = "  +1.51%  "

# Row 17 - ShibaInu
Set-PriceText "D17" "0.00001131"
$ws.Range("E17").Value = "  +1.97%  "

# Row 18 - TRON
Set-PriceText "D18" "0.06713"
$ws.Range("E18").Value = "  +0.18%  "

# Row 19 - Dai
Set-PriceText "D19" "0.9976"
$ws.Range("E19").Value = "  -0.17%  "

# Row 20 - Litecoin
Set-PriceText "D20" "84.39"
$ws.Range("E20").Value = "  +3.98%  "

# Row 21 - Avalanche
Set-PriceText "D21" "17.28"
$ws.Range("E21").Value = "  +4.65%  "

# Row 22 - Uniswap
Set-PriceText "D22" "6.386"
$ws.Range("E22").Value = "  +4.24%  "

# Row 23 - Cosmos (price unchanged, volume refreshed)
$ws.Range("E23").Value = "  +7.09%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "24.912.18"
$ws.Range("E24").Value = "  +2.07%  "

# Row 25 - Toncoin (price unchanged, volume refreshed)
$ws.Range("E25").Value = "  -0.33%  "

# Row 26 - LidoDAOToken (price unchanged, volume refreshed)
$ws.Range("E26").Value = "  +4.56%  "

# Row 27 - EthereumClassic
Set-PriceText "D27" "20.43"
$ws.Range("E27").Value = "  +4.54%  "

# Row 28 - Monero
Set-PriceText "D28" "149.99"
$ws.Range("E28").Value = "  -1.80%  "

# Row 29 - BitcoinCash
Set-PriceText "D29" "133.43"
$ws.Range("E29").Value = "  +5.03%  "

# Row 30 - WrappedliquidstakedEther2.0 (volume unchanged)
$ws.Range("D30").Value = "1.894.81"

# Row 31 - ImmutableX
Set-PriceText "D31" "1.231"
$ws.Range("E31").Value = "  +25.81%  "

# Row 32 - Filecoin
Set-PriceText "D32" "6.859"
$ws.Range("E32").Value = "  +8.03%  "

# Row 33 - HuobiToken (price unchanged, volume refreshed)
$ws.Range("E33").Value = "  +4.40%  "

# Row 34 - Aptos (price unchanged, volume refreshed)
$ws.Range("E34").Value = "  +11.37%  "

# Row 35 / 36 - coin ranking order swapped: Stellar now ranks above WEMIXTOKEN.
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-PriceText "D35" "0.08800"
$ws.Range("E35").Value = "  +4.09%  "

$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-PriceText "D36" "1.773"
$ws.Range("E36").Value = "  +4.25%  "

# Row 37 - InternetComputer(DFINITY)
Set-PriceText "D37" "5.637"
$ws.Range("E37").Value = "  +4.53%  "

# Row 38 - Hedera
Set-PriceText "D38" "0.06658"
$ws.Range("E38").Value = "  +2.01%  "

# Row 39 - FraxShare
Set-PriceText "D39" "9.174"
$ws.Range("E39").Value = "  +2.73%  "

# Row 40 - VeChain
Set-PriceText "D40" "0.02416"
$ws.Range("E40").Value = "  +3.25%  "

# Row 41 - Algorand
Set-PriceText "D41" "0.2230"
$ws.Range("E41").Value = "  +5.19%  "

# Row 42 - TrustWalletToken
Set-PriceText "D42" "1.283"
$ws.Range("E42").Value = "  +1.68%  "

# Row 43 - TheSandbox (price unchanged, volume refreshed)
$ws.Range("E43").Value = "  +4.55%  "

# Row 44 - Frax
Set-PriceText "D44" "0.9981"
$ws.Range("E44").Value = "  -0.10%  "

# Row 45 - EnergySwap
Set-PriceText "D45" "13.92"
$ws.Range("E45").Value = "  +5.24%  "

# Row 46 - Decentraland
Set-PriceText "D46" "0.6159"
$ws.Range("E46").Value = "  +3.23%  "

# Row 47 - PancakeSwap
Set-PriceText "D47" "3.823"
$ws.Range("E47").Value = "  +1.26%  "

# Row 48 - NEARProtocol
Set-PriceText "D48" "2.125"
$ws.Range("E48").Value = "  +4.65%  "

# Row 49 - Quant
Set-PriceText "D49" "129.83"
$ws.Range("E49").Value = "  +1.79%  "

# Row 50 - Cronos
Set-PriceText "D50" "0.07313"
$ws.Range("E50").Value = "  +1.74%  "

# Row 51 - Aave
Set-PriceText "D51" "80.01"
$ws.Range("E51").Value = "  +5.24%  "
